# CdCity.xlsx maintenance edit:
# Remove the "JcicCityCode" (聯徵用縣市代碼) field definition row from the
# DBD layout sheet. Deleting the whole worksheet row shifts every row
# below it up by one and lets Excel recompute the shared-string table,
# which is exactly what the target revision shows (row 22 removed, rows
# 23-36 become rows 22-35).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# Row 22 holds: SEQ 15 | JcicCityCode | 聯徵用縣市代碼 | VARCHAR2 | (blank) | (blank) | 刪除
# Delete it entirely -- cells/rows below ripple upward automatically.
$ws.Rows(22).Delete()

# Match the author's final on-screen selection/scroll position.
$ws.Range("A9").Select()
$ws.Range("B16").Select()
